$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking strings (e.g. "1.00", "3.50")
# keep their exact text representation instead of being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '72.002.71'
$ws.Range("E2").Value = '  +0.61%  '

$ws.Range("D3").Value = '4.006.48'
$ws.Range("E3").Value = '  +0.11%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = '533.34'
$ws.Range("E5").Value = '  +0.93%  '

$ws.Range("D6").Value = '152.35'
$ws.Range("E6").Value = '  +2.46%  '

$ws.Range("D7").Value = '0.696'
$ws.Range("E7").Value = '  +11.35%  '

$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("D9").Value = '0.747'
$ws.Range("E9").Value = '  +1.52%  '

$ws.Range("D10").Value = '0.171'
$ws.Range("E10").Value = '  -3.06%  '

$ws.Range("D11").Value = '0.0000325'
$ws.Range("E11").Value = '  -5.72%  '

$ws.Range("D12").Value = '47.66'
$ws.Range("E12").Value = '  +7.24%  '

$ws.Range("D13").Value = '10.64'
$ws.Range("E13").Value = '  -0.48%  '

$ws.Range("D14").Value = '4.641.97'
$ws.Range("E14").Value = '  -0.08%  '

$ws.Range("D15").Value = '3.986.26'
$ws.Range("E15").Value = '  -1.04%  '

$ws.Range("D16").Value = '13.92'
$ws.Range("E16").Value = '  -2.26%  '

$ws.Range("D17").Value = '20.42'
$ws.Range("E17").Value = '  -4.29%  '

$ws.Range("E18").Value = '  -1.11%  '

$ws.Range("D19").Value = '1.18'
$ws.Range("E19").Value = '  -2.70%  '

$ws.Range("D20").Value = '71.878.13'
$ws.Range("E20").Value = '  +0.59%  '

$ws.Range("D21").Value = '426.45'
$ws.Range("E21").Value = '  -3.54%  '

$ws.Range("D22").Value = '97.91'
$ws.Range("E22").Value = '  +4.34%  '

$ws.Range("D23").Value = '3.50'
$ws.Range("E23").Value = '  -2.45%  '

$ws.Range("D24").Value = '4.18'
$ws.Range("E24").Value = '  +2.00%  '

$ws.Range("D25").Value = '14.33'
$ws.Range("E25").Value = '  -0.03%  '

$ws.Range("D26").Value = '11.10'
$ws.Range("E26").Value = '  -9.50%  '

$ws.Range("D27").Value = '10.65'
$ws.Range("E27").Value = '  -2.25%  '

$ws.Range("E28").Value = '  +1.47%  '

$ws.Range("D29").Value = '36.68'
$ws.Range("E29").Value = '  -0.74%  '

$ws.Range("D30").Value = '3.59'
$ws.Range("E30").Value = '  +24.03%  '

$ws.Range("D31").Value = '13.35'
$ws.Range("E31").Value = '  -1.88%  '

$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = '0.129'
$ws.Range("E32").Value = '  -0.39%  '

$ws.Range("B33").Value = 'Bittensor'
$ws.Range("C33").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D33").Value = '674.04'
$ws.Range("E33").Value = '  -3.84%  '

$ws.Range("D34").Value = '7.07'
$ws.Range("E34").Value = '  +1.42%  '

$ws.Range("D35").Value = '65.66'
$ws.Range("E35").Value = '  -2.06%  '

$ws.Range("D36").Value = '42.73'
$ws.Range("E36").Value = '  +4.09%  '

$ws.Range("D37").Value = '0.426'
$ws.Range("E37").Value = '  -4.10%  '

$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").Value = '0.152'
$ws.Range("E38").Value = '  +0.06%  '

$ws.Range("B39").Value = 'PEPE'
$ws.Range("C39").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D39").Value = ('0.0' + [string][char]0x2083 + '0824')
$ws.Range("E39").Value = '  -9.71%  '

$ws.Range("D40").Value = '3.45'
$ws.Range("E40").Value = '  -4.75%  '

$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  -0.10%  '

$ws.Range("B42").Value = 'WEMIXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D42").Value = '3.28'
$ws.Range("E42").Value = '  +4.64%  '

$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").Value = '0.998'
$ws.Range("E43").Value = '  -0.20%  '

$ws.Range("D44").Value = '0.0487'
$ws.Range("E44").Value = '  -1.65%  '

$ws.Range("E45").Value = '  +2.91%  '

$ws.Range("D46").Value = '9.66'
$ws.Range("E46").Value = '  +4.20%  '

$ws.Range("D47").Value = '2.61'
$ws.Range("E47").Value = '  -10.23%  '

$ws.Range("D48").Value = '3.36'
$ws.Range("E48").Value = '  -5.42%  '

$ws.Range("D49").Value = '2.98'
$ws.Range("E49").Value = '  -7.87%  '

$ws.Range("B50").Value = 'FLOKI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D50").Value = '0.000270'
$ws.Range("E50").Value = '  -5.35%  '

$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").Value = '144.17'
$ws.Range("E51").Value = '  +1.04%  '

# Restore original (default) cell style on column D now that text values are set,
# so no lingering "Text" number-format style is left applied to the cells.
$ws.Range("D2:D51").ClearFormats()
